$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated extrapolation calibration values (removing sub-$5 price noise)

$ws.Range("D3").Value = 114681.8544033737
$ws.Range("E3").Value = 0.04850074704636939
$ws.Range("F3").Value = 0.2224122072687891
$ws.Range("G3").Value = -2.571228540895307
$ws.Range("H3").Value = 33.2178303362922

$ws.Range("D5").Value = 116108.8641132673
$ws.Range("E5").Value = 0.01719121215503349
$ws.Range("F5").Value = 0.2401399839392263
$ws.Range("G5").Value = -1.265845000868734
$ws.Range("H5").Value = 12.5086040286615

$ws.Range("D6").Value = 116748.8310009253
$ws.Range("E6").Value = 0.004653008165384166
$ws.Range("F6").Value = 0.2737839248593684
$ws.Range("G6").Value = -1.973925232086058
$ws.Range("H6").Value = 17.69395457043268

$ws.Range("D7").Value = 117266.6974678525
$ws.Range("E7").Value = -0.002184313165384341
$ws.Range("F7").Value = 0.267374056671252
$ws.Range("G7").Value = -1.559151150081153
$ws.Range("H7").Value = 14.23413940883779

$ws.Range("D8").Value = 118833.9973636744
$ws.Range("E8").Value = -0.02669963728240702
$ws.Range("F8").Value = 0.2261499164594412
$ws.Range("G8").Value = -0.941033161531279
$ws.Range("H8").Value = 7.125176421378931

$ws.Range("D9").Value = 120551.0117654942
$ws.Range("E9").Value = -0.06884031825034642
$ws.Range("F9").Value = 0.3587020288023273
$ws.Range("G9").Value = -1.358103254148493
$ws.Range("H9").Value = 8.899885722194352

$ws.Range("D10").Value = 121869.5443124926
$ws.Range("E10").Value = -0.1120817811249724
$ws.Range("F10").Value = 0.4691483585477071
$ws.Range("G10").Value = -1.71764757331361
$ws.Range("H10").Value = 8.735833435738616

$ws.Range("D14").Value = 113875.9553570923
$ws.Range("E14").Value = 0.1825808659858967
$ws.Range("F14").Value = 0.1508774629066905
$ws.Range("G14").Value = -0.2435196048899383
$ws.Range("H14").Value = 6.600322888226663

$ws.Range("D17").Value = 113605.0301608585
$ws.Range("E17").Value = 0.1406442854260334
$ws.Range("F17").Value = 0.1350245255001157
$ws.Range("G17").Value = -0.6921808950884332
$ws.Range("H17").Value = 6.475307955945548

$ws.Range("D18").Value = 114723.6678293223
$ws.Range("E18").Value = 0.1087880492824932
$ws.Range("F18").Value = 0.1726661219483923
$ws.Range("G18").Value = -0.1988639614546169
$ws.Range("H18").Value = 4.83555270284583
